$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from SNOMED CT[ N]" sheets to "Include #<n>" (0-based)
for ($i = 0; $i -le 8; $i++) {
    $oldName = "Include from SNOMED CT"
    if ($i -gt 0) { $oldName = "Include from SNOMED CT " + ($i + 1) }
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = "Include #$i"
}

# 2. Update Metadata sheet values
$meta = $wb.Worksheets.Item("Metadata")

# URL: pythia -> cicada
$meta.Range("B2").Value2 = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/vaccine-medication-codes-snomed"

# Date: updated timestamp
$meta.Range("B8").Value2 = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" row between "Contact" (row 10) and "Description" (row 11)
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$meta.Range("A11").Value2 = "Jurisdiction"
$meta.Range("B11").Value2 = ""

Write-Host "Edit complete"
